$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.374.08"
$ws.Range("E2").Value = "  -2.37%  "
$ws.Range("D3").Value = "2.294.82"
$ws.Range("E3").Value = "  -1.79%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'496.45"
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("D6").Value = "'127.43"
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.529"
$ws.Range("E8").Value = "  -1.24%  "
$ws.Range("D9").Value = "2.294.06"
$ws.Range("E9").Value = "  -2.08%  "
$ws.Range("D10").Value = "'0.0944"
$ws.Range("E10").Value = "  -3.37%  "
$ws.Range("D12").Value = "'0.322"
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("D13").Value = "'4.64"
$ws.Range("E13").Value = "  -2.90%  "
$ws.Range("D14").Value = "2.695.00"
$ws.Range("E14").Value = "  -2.19%  "
$ws.Range("D15").Value = "'21.68"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "54.274.49"
$ws.Range("E16").Value = "  -2.52%  "
$ws.Range("D17").Value = "'0.0000130"
$ws.Range("E17").Value = "  -1.58%  "
$ws.Range("D18").Value = "2.289.43"
$ws.Range("E18").Value = "  -2.40%  "
$ws.Range("D19").Value = "'9.95"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").Value = "'4.06"
$ws.Range("E20").Value = "  +1.86%  "
$ws.Range("D21").Value = "'303.50"
$ws.Range("E21").Value = "  -1.93%  "
$ws.Range("D22").Value = "'6.39"
$ws.Range("E22").Value = "  +3.28%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'63.90"
$ws.Range("E24").Value = "  -2.04%  "
$ws.Range("D26").Value = "'0.375"
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("D27").Value = "2.367.81"
$ws.Range("E27").Value = "  -3.42%  "
$ws.Range("D28").Value = "'0.149"
$ws.Range("E28").Value = "  +2.27%  "
$ws.Range("D29").Value = "'7.16"
$ws.Range("E29").Value = "  +1.17%  "
$ws.Range("D30").Value = "'166.01"
$ws.Range("E30").Value = "  -3.35%  "
$ws.Range("E31").Value = "  -1.87%  "
$ws.Range("D32").Value = "0.0₃0683"
$ws.Range("E32").Value = "  -2.44%  "
$ws.Range("D33").Value = "'5.88"
$ws.Range("E33").Value = "  +1.85%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("E36").Value = "  +1.58%  "
$ws.Range("D37").Value = "'17.63"
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("D38").Value = "'1.19"
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("D39").Value = "'0.878"
$ws.Range("E39").Value = "  +7.33%  "
$ws.Range("D40").Value = "'3.65"
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("D41").Value = "'35.45"
$ws.Range("E41").Value = "  -1.73%  "
$ws.Range("D42").Value = "'0.376"
$ws.Range("E42").Value = "  +1.46%  "
$ws.Range("E43").Value = "  +1.93%  "
$ws.Range("D44").Value = "'3.36"
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'125.98"
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'4.83"
$ws.Range("E46").Value = "  +3.01%  "
$ws.Range("D47").Value = "'0.0891"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").Value = "'0.547"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("D49").Value = "'238.22"
$ws.Range("E49").Value = "  +0.87%  "
$ws.Range("D50").Value = "'0.0481"
$ws.Range("E50").Value = "  +1.41%  "
$ws.Range("D51").Value = "'0.0205"
$ws.Range("E51").Value = "  +0.19%  "
